$d = $word.ActiveDocument

# Locate the SourceCode paragraph holding the R console output that
# follows "head(painters)" (the "##  ... " table dump) and remove the
# whole paragraph, including its paragraph mark, so the document flows
# directly from "head(painters)" into "Question 1:".
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t.StartsWith("##               Composition Drawing Colour Expression School")) {
        $p.Range.Delete()
        break
    }
}
